$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new Gherkin test-data row (row 6) ---------------------------------
$ws.Range("A6").Value = "mail.unregistered@test.com"
$ws.Range("B6").Value = "ABCde!12350"
$ws.Range("C6").Value = "incorrect account"

# Give A6 the same "hyperlink" cell formatting that A2 already has, and wire
# up a live mailto: hyperlink on it (mirrors the A2 / rId1 pattern).
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:mail.unregistered@test.com") | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Resize columns A and C to comfortably fit the new values --------------
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15

# --- Move the selection like the author's last recorded cursor position ----
$ws.Range("C7").Select() | Out-Null
